$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"0.0001021024915524027"
$ws.Range("C2").Value = [double]"1.829203455372408e-12"
$ws.Range("D2").Value = [double]"0.7527432677738641"
$ws.Range("E2").Value = [double]"0.4942365360607697"
$ws.Range("G2").Value = [double]"1.247081906328015"

$ws.Range("B3").Value = [double]"0.2917716402565462"
$ws.Range("C3").Value = [double]"0.306821227259698"
$ws.Range("D3").Value = [double]"0.7527432677738641"
$ws.Range("E3").Value = [double]"0.4942365360607697"
$ws.Range("G3").Value = [double]"1.845572671350878"

$ws.Range("B4").Value = [double]"0.6606524410359556"
$ws.Range("C4").Value = [double]"0.306821227259698"
$ws.Range("D4").Value = [double]"0.1494219747398047"
$ws.Range("E4").Value = [double]"0.4942365360607697"
$ws.Range("G4").Value = [double]"1.611132179096228"
